# fixed checking reminder_date is empty
#
# The "reminder_date" template sample in row 3 (H3) is cleared so the
# workbook's consuming script can be exercised against an empty
# reminder_date. While touching the sample data, the bcc_emails column
# (C2/C3) is emptied and the to_emails address used on row 2 (A2) is
# switched to the same address already used on row 3 (A3),
# nmhillusion@hotmail.com, retargeting the mail-to hyperlinks accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- to_emails (A2 / A3): point both rows at nmhillusion@hotmail.com ---

# Drop every existing hyperlink in the sheet (A2, C2, C3, A3) so we can
# rebuild only the ones that should remain.
$ws.Range("A1").Hyperlinks.Delete()

$ws.Range("A2").Value = "nmhillusion@hotmail.com"

$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:nmhillusion@hotmail.com")
$ws.Range("A2").NumberFormat = "@"

$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:nmhillusion@hotmail.com")
$ws.Range("A3").NumberFormat = "@"

# --- bcc_emails (C2 / C3): clear the sample addresses ---
$ws.Range("C2").ClearContents()
$ws.Range("C3").ClearContents()

# --- reminder_date (H3): clear entirely so the empty-check can be tested ---
$ws.Range("H3").Clear()

# Matches the saved selection recorded in the workbook after the edit.
[void]$ws.Range("C6").Select()
